# Corrected Cap Part Number on BOM
#
# Row 5 of the DueShield BOM sheet (designators C7-C13, footprint 0805-CAP,
# qty 7) listed the wrong part for the 0.01uF/10nF decoupling capacitor.
# Replace the Mfg Part Number / Source / Description with the corrected
# values and highlight the corrected row in yellow so it stands out.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophes force these in as literal text (keeping the existing
# "quote prefix" text formatting already used throughout this BOM column)
# rather than creating a plain numeric/general style.
$ws.Range("A5").Value = "'581-08055C103J"
$ws.Range("B5").Value = "'MOUSER"
$ws.Range("C5").Value = "'CAPACITOR, 0.01uF, 50V"

# Highlight the corrected row in yellow.
$ws.Range("A5:C5").Interior.Color = 65535

# Leave the selection where the author left it when saving.
$ws.Range("C10").Select()
